$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: Walsh & Patterson (Ottawa precipitation/temperature trends) ---
$ws.Range("A37").Value = "Precipitation and Temperature Trends and Cycles Derived from Historical 1890-2019 Weather Data for the City of Ottawa, Ontario, Canada"
$ws.Range("B37").Value = "Environments"
$ws.Range("C37").Value = 2022
$ws.Range("D37").Value = "Walsh & Patterson"
$ws.Range("F37").Value = "NA"
$ws.Range("G37").Value = "Cimatic trend in the last century in Ottawa"
$ws.Range("I37").Value = "yes"
$ws.Range("J37").Value = "yes"
$ws.Range("K37").Value = "yes"

# --- Row 38: Zhai et al. (climate projections Ottawa) ---
$ws.Range("A38").Value = "Future projections of temperature changes in Ottawa, Canada through stepwise clustered downscaling of multiple GCMs under RCPs"
$ws.Range("B38").Value = "Climate Dynamics"
$ws.Range("C38").Value = 2019
$ws.Range("D38").Value = "Zhai, Huang, Wang, Zhou, Lu, Li"
$ws.Range("F38").Value = "NA"
$ws.Range("G38").Value = "Climate projections for Ottawa"
$ws.Range("I38").Value = "yes"
$ws.Range("J38").Value = "yes"
$ws.Range("K38").Value = "yes"

# --- Row 39: Geiser (hibernation definition) ---
$ws.Range("A39").Value = "Hibernation"
$ws.Range("B39").Value = "Current biology"
$ws.Range("C39").Value = 2013
$ws.Range("D39").Value = "Geiser"
$ws.Range("F39").Value = "NA"
$ws.Range("G39").Value = "Hibernation' definition"
$ws.Range("I39").Value = "yes"
$ws.Range("J39").Value = "yes"
$ws.Range("K39").Value = "yes"

# --- Row 40: Nedergaard & Cannon (mammalian hibernation) ---
$ws.Range("A40").Value = "Mammalian hibernation"
$ws.Range("B40").Value = "Phil. Trans. R. Soc. Lond. B"
$ws.Range("C40").Value = 1990
$ws.Range("D40").Value = "Nedergaard & Cannon"
$ws.Range("F40").Value = "NA"
$ws.Range("G40").Value = "Mammalian hibernation"
$ws.Range("I40").Value = "yes"
$ws.Range("J40").Value = "yes"
$ws.Range("K40").Value = "yes"

# --- Row 41: Durant et al. (climate/predator phenology mismatch) ---
$ws.Range("A41").Value = "Climate and the match or mismatch between predator requirements and resource availability"
$ws.Range("D41").Value = "Durant, Hjermann, Ottersen, Stenseth"
$ws.Range("B41").Value = "Clim. Res."
$ws.Range("C41").Value = 2007
$ws.Range("F41").Value = "no"
$ws.Range("G41").Value = "changes in phenology as a universal response to climate change"
$ws.Range("I41").Value = "yes"
$ws.Range("J41").Value = "yes"
$ws.Range("K41").Value = "yes"

# --- Fix up cell formatting (styles) to match the target workbook look ---
# F column ("Read" -> NA) takes the "filled" style used elsewhere (copy from F36)
$ws.Range("F36").Copy()
$ws.Range("F37").PasteSpecial(-4122)
$ws.Range("F38").PasteSpecial(-4122)
$ws.Range("F39").PasteSpecial(-4122)
$ws.Range("F40").PasteSpecial(-4122)

# I/J/K columns ("yes") take the "filled" style (copy from I36/J36/K36)
$ws.Range("I36").Copy()
$ws.Range("I37").PasteSpecial(-4122)
$ws.Range("I39").PasteSpecial(-4122)
$ws.Range("I40").PasteSpecial(-4122)
$ws.Range("I41").PasteSpecial(-4122)

$ws.Range("J36").Copy()
$ws.Range("J37").PasteSpecial(-4122)
$ws.Range("J39").PasteSpecial(-4122)

$ws.Range("K36").Copy()
$ws.Range("K37").PasteSpecial(-4122)
$ws.Range("K39").PasteSpecial(-4122)
$ws.Range("K40").PasteSpecial(-4122)
$ws.Range("K41").PasteSpecial(-4122)

# D39 (newly used cell) picks up the standard author-column style
$ws.Range("D36").Copy()
$ws.Range("D39").PasteSpecial(-4122)

# B41 (Journal column) takes the alternate "bold-ish" style used by B31
$ws.Range("B31").Copy()
$ws.Range("B41").PasteSpecial(-4122)

# G41 (new cell, Point of the paper) takes the standard style
$ws.Range("G38").Copy()
$ws.Range("G41").PasteSpecial(-4122)

# G42 loses its former "centered bold" style now that row 41 above it is filled in
$ws.Range("G38").Copy()
$ws.Range("G42").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Scroll position / selection, matching where review left off ---
$ws.Range("B12").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G45").Select()
